$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new broker/quarter data rows (68-76) in columns A, B, F.
# (Columns C/D are left blank for these rows, matching the source edit.)

$ws.Range("A68").Value = "SHS"
$ws.Range("B68").Value = "2Q25"
$ws.Range("F68").Value = 462.96722006800002
$ws.Range("F68").NumberFormat = "#,##0"

$ws.Range("A69").Value = "VIX"
$ws.Range("B69").Value = "2Q25"
$ws.Range("F69").Value = 1602.5073289380002
$ws.Range("F69").NumberFormat = "#,##0"

$ws.Range("A70").Value = "SHS"
$ws.Range("B70").Value = "4Q24"
$ws.Range("F70").Value = 281.90273591299996
$ws.Range("F70").NumberFormat = "#,##0"

$ws.Range("A71").Value = "SHS"
$ws.Range("B71").Value = "1Q25"
$ws.Range("F71").Value = 325.68608381000001
$ws.Range("F71").NumberFormat = "#,##0"

$ws.Range("A72").Value = "VIX"
$ws.Range("B72").Value = "4Q24"
$ws.Range("F72").Value = 133.6666776649999
$ws.Range("F72").NumberFormat = "#,##0"

$ws.Range("A73").Value = "VIX"
$ws.Range("B73").Value = "1Q25"
$ws.Range("F73").Value = 465.20486131000007
$ws.Range("F73").NumberFormat = "#,##0"

$ws.Range("A74").Value = "CTS"
$ws.Range("B74").Value = "4Q24"
$ws.Range("F74").Value = 73.276574532000055
$ws.Range("F74").NumberFormat = "#,##0"

$ws.Range("A75").Value = "CTS"
$ws.Range("B75").Value = "1Q25"
$ws.Range("F75").Value = 132.5779621770001
$ws.Range("F75").NumberFormat = "#,##0"

$ws.Range("A76").Value = "CTS"
$ws.Range("B76").Value = "2Q25"
$ws.Range("F76").Value = 217.38459526100004
$ws.Range("F76").NumberFormat = "#,##0"

# Update the view so the new row F68 is the active selection (matches the
# saved workbook's cursor position after the data entry).
$ws.Range("F68").Select() | Out-Null
